$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G header
$ws.Range("G1").Value = "results"

# Rows where the analysis actually ran (CMH test rows) get an #N/A error
# in the new "results" column; rows where "Analysis Did Not Run" get the
# same text repeated in column G.
$naRows = @(2,3,4,5,6,7,11,12,13,15,16,17)
foreach ($r in $naRows) {
    $ws.Cells.Item($r, 7).Value = "#N/A"
}

$notRunRows = @(8,9,10,14,18)
foreach ($r in $notRunRows) {
    $ws.Cells.Item($r, 7).Value = "Analysis Did Not Run"
}
